$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.569.52'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('D3').Value = '3.065.53'
$ws.Range('E3').Value = '  +2.83%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '386.50'
$ws.Range('E5').Value = '  +1.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '103.16'
$ws.Range('E6').Value = '  -0.13%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.544'
$ws.Range('E7').Value = '  -0.31%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.585'
$ws.Range('E9').Value = '  -1.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.84'
$ws.Range('E10').Value = '  +0.45%  '
$ws.Range('E11').Value = '  +0.30%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0862'
$ws.Range('E12').Value = '  +0.31%  '
$ws.Range('D13').Value = '3.549.80'
$ws.Range('E13').Value = '  +2.67%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.66'
$ws.Range('E14').Value = '  +1.25%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.78'
$ws.Range('E15').Value = '  +0.01%  '
$ws.Range('D16').Value = '3.056.48'
$ws.Range('E16').Value = '  +2.30%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.975'
$ws.Range('E17').Value = '  -2.08%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.73'
$ws.Range('E18').Value = '  -3.33%  '
$ws.Range('D19').Value = '51.621.50'
$ws.Range('E19').Value = '  +0.41%  '
$ws.Range('E20').Value = '  +2.68%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.47'
$ws.Range('E21').Value = '  -0.92%  '
$ws.Range('D22').Value = '0.0₃0969'
$ws.Range('E22').Value = '  +0.84%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.18'
$ws.Range('E23').Value = '  -0.41%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '268.14'
$ws.Range('E24').Value = '  +0.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.16'
$ws.Range('E25').Value = '  -1.43%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.21'
$ws.Range('E26').Value = '  +4.85%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '26.91'
$ws.Range('E27').Value = '  +3.45%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.31'
$ws.Range('E28').Value = '  -1.17%  '
$ws.Range('E29').Value = '  +2.04%  '
$ws.Range('E30').Value = '  +0.17%  '
$ws.Range('E31').Value = '  -1.55%  '
$ws.Range('E32').Value = '  -0.22%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '34.82'
$ws.Range('E33').Value = '  +0.55%  '
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '50.03'
$ws.Range('E35').Value = '  -2.94%  '
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0449'
$ws.Range('E36').Value = '  +2.28%  '
$ws.Range('E37').Value = '  -0.16%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.33'
$ws.Range('E38').Value = '  +2.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.292'
$ws.Range('E39').Value = '  +8.63%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.87'
$ws.Range('E40').Value = '  +1.92%  '
$ws.Range('B41').Value = 'Celestia'
$ws.Range('C41').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '16.92'
$ws.Range('E41').Value = '  +0.91%  '
$ws.Range('E42').Value = '  +0.60%  '
$ws.Range('E43').Value = '  -0.44%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '125.35'
$ws.Range('E44').Value = '  +0.57%  '
$ws.Range('E45').Value = '  +1.89%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '22.07'
$ws.Range('E46').Value = '  +2.77%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.09'
$ws.Range('E47').Value = '  +2.51%  '
$ws.Range('E48').Value = '  +1.15%  '
$ws.Range('D49').Value = '2.037.15'
$ws.Range('E49').Value = '  +0.11%  '
$ws.Range('D50').Value = '3.363.83'
$ws.Range('E50').Value = '  +2.27%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.208'
$ws.Range('E51').Value = '  +7.68%  '
